$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the target paragraph ("Github does have a so-called abuse
#    prevention mechanism ... worth mentioning.")
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*so-called abuse prevention mechanism*") {
        $target = $p
        break
    }
}

# ------------------------------------------------------------------
# 2. Split "It's not obvious this matters, but worth mentioning." so
#    "this matters" sits in its own run (grammar-check highlight).
# ------------------------------------------------------------------
$full = $target.Range.Text
$idx = $full.IndexOf("this matters")
$pstart = $target.Range.Start
$subStart = $pstart + $idx
$subEnd = $subStart + ("this matters").Length
$sub = $d.Range($subStart, $subEnd)
$sub.Bold = 1
$sub.Bold = 0

# ------------------------------------------------------------------
# 3. Insert two new list paragraphs after the target paragraph, same
#    list style / level / numbering as the target.
# ------------------------------------------------------------------
$target.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($target.Index + 1)
$p2.Range.InsertAfter("Scraped emails were double-cleaned using emaillistverify and mailgun.")

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($p2.Index + 1)
$p3.Range.InsertAfter("Undefined scrape date is before march 8")

Write-Output "done"
